$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 10419606
$ws.Range("I33").Value = 17857778
$ws.Range("J33").Value = 6164.7
$ws.Range("K33").Value = 17857778
$ws.Range("L33").Value = 6164.7
$ws.Range("M33").Value = -17857549
$ws.Range("N33").Value = -6622.7
$ws.Range("H43").Value = 6083.4546
$ws.Range("J43").Value = 8500.5
$ws.Range("L43").Value = 8500.5
$ws.Range("N43").Value = -8638.5
$ws.Range("H99").Value = 1198.375
$ws.Range("I99").Value = 1198.375
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3595.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2097.125
$ws.Range("N99").ClearContents()
$ws.Range("H116").Value = 2803.8513
$ws.Range("I116").Value = 2725.2957
$ws.Range("J116").Value = 4663
$ws.Range("K116").Value = 2725.2957
$ws.Range("L116").Value = 4663
$ws.Range("M116").Value = 716.7042999999999
$ws.Range("N116").Value = -11547
$ws.Range("H134").Value = 83833.336
$ws.Range("J134").Value = 83833.336
$ws.Range("L134").Value = 83833.336
$ws.Range("N134").Value = -93973.336

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3362.182
$ws.Range("I2").Value = 968.3333
$ws.Range("J2").Value = 6234.8
$ws.Range("K2").Value = 968.3333
$ws.Range("L2").Value = 6234.8
$ws.Range("M2").Value = -855.3333
$ws.Range("N2").Value = -6460.8
$ws.Range("H32").Value = 14971.948
$ws.Range("I32").Value = 5158.75
$ws.Range("K32").Value = 5158.75
$ws.Range("M32").Value = -4871.75
$ws.Range("H102").Value = 15550.235
$ws.Range("I102").Value = 4133.25
$ws.Range("K102").Value = 4133.25
$ws.Range("M102").Value = -2511.25
$ws.Range("H116").Value = 3362.182
$ws.Range("I116").Value = 968.3333
$ws.Range("J116").Value = 6234.8
$ws.Range("K116").Value = 968.3333
$ws.Range("L116").Value = 6234.8
$ws.Range("M116").Value = 1325.6667
$ws.Range("N116").Value = -10822.8
$ws.Range("H134").Value = 96000
$ws.Range("J134").Value = 96000
$ws.Range("L134").Value = 96000
$ws.Range("N134").Value = -106140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3362.182
$ws.Range("I3").Value = 968.3333
$ws.Range("J3").Value = 6234.8
$ws.Range("K3").Value = 968.3333
$ws.Range("L3").Value = 6234.8
$ws.Range("M3").Value = -854.3333
$ws.Range("N3").Value = -6462.8
$ws.Range("H20").Value = 22342.062
$ws.Range("I20").Value = 10562.885
$ws.Range("K20").Value = 10562.885
$ws.Range("M20").Value = -10315.885
$ws.Range("H99").Value = 1611.8823
$ws.Range("I99").Value = 1171
$ws.Range("K99").Value = 1171
$ws.Range("M99").Value = 327
$ws.Range("H134").Value = 15438.73
$ws.Range("I134").Value = 8647.8125
$ws.Range("J134").Value = 26304.2
$ws.Range("K134").Value = 25943.4375
$ws.Range("L134").Value = 78912.60000000001
$ws.Range("M134").Value = -23408.4375
$ws.Range("N134").Value = -83982.60000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 357.8889
$ws.Range("J2").Value = 479.8
$ws.Range("L2").Value = 479.8
$ws.Range("N2").Value = -705.8
$ws.Range("H14").Value = 520.6
$ws.Range("J14").Value = 475
$ws.Range("L14").Value = 475
$ws.Range("N14").Value = -815
$ws.Range("H15").Value = 30225
$ws.Range("I15").Value = 450
$ws.Range("K15").Value = 450
$ws.Range("M15").Value = -280
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H99").Value = 7941.4116
$ws.Range("I99").Value = 2409.1667
$ws.Range("J99").Value = 10959
$ws.Range("K99").Value = 2409.1667
$ws.Range("L99").Value = 10959
$ws.Range("M99").Value = -911.1667000000002
$ws.Range("N99").Value = -13955
$ws.Range("H126").Value = 7941.4116
$ws.Range("I126").Value = 2409.1667
$ws.Range("J126").Value = 10959
$ws.Range("K126").Value = 7227.500100000001
$ws.Range("L126").Value = 32877
$ws.Range("M126").Value = -4757.500100000001
$ws.Range("N126").Value = -37817

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 91093180
$ws.Range("I4").Value = 143003020
$ws.Range("J4").Value = 250960
$ws.Range("K4").Value = 429009060
$ws.Range("L4").Value = 752880
$ws.Range("M4").Value = -429008948
$ws.Range("N4").Value = -753104
$ws.Range("H11").Value = 1603.75
$ws.Range("I11").Value = 1788.5
$ws.Range("J11").Value = 1049.5
$ws.Range("K11").Value = 5365.5
$ws.Range("L11").Value = 3148.5
$ws.Range("M11").Value = -5225.5
$ws.Range("N11").Value = -3428.5
$ws.Range("H46").Value = 715
$ws.Range("I46").Value = 591
$ws.Range("J46").Value = 1025
$ws.Range("K46").Value = 1773
$ws.Range("L46").Value = 3075
$ws.Range("M46").Value = -1682
$ws.Range("N46").Value = -3257
$ws.Range("H56").Value = 200005340
$ws.Range("I56").Value = 200005340
$ws.Range("K56").Value = 200005340
$ws.Range("M56").Value = -200004810
$ws.Range("H128").Value = 203778.6
$ws.Range("I128").Value = 203778.6
$ws.Range("K128").Value = 611335.8
$ws.Range("M128").Value = -606355.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2250645
$ws.Range("I14").Value = 3000563.2
$ws.Range("K14").Value = 3000563.2
$ws.Range("M14").Value = -3000395.2
$ws.Range("H17").Value = 1579.9333
$ws.Range("J17").Value = 1671.3572
$ws.Range("L17").Value = 1671.3572
$ws.Range("N17").Value = -2007.3572
$ws.Range("H23").Value = 430.27777
$ws.Range("J23").Value = 1083.3334
$ws.Range("L23").Value = 1083.3334
$ws.Range("N23").Value = -1529.3334
$ws.Range("H102").Value = 7097.3184
$ws.Range("I102").Value = 4918.1055
$ws.Range("K102").Value = 4918.1055
$ws.Range("M102").Value = -3296.1055
$ws.Range("H135").Value = 155058.77
$ws.Range("J135").Value = 155058.77
$ws.Range("L135").Value = 155058.77
$ws.Range("N135").Value = -165198.77

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4365.485
$ws.Range("I22").Value = 2144.158
$ws.Range("K22").Value = 2144.158
$ws.Range("M22").Value = -1849.158
$ws.Range("H27").Value = 4365.485
$ws.Range("I27").Value = 2144.158
$ws.Range("K27").Value = 2144.158
$ws.Range("M27").Value = -2037.158
$ws.Range("H46").Value = 190908.78
$ws.Range("J46").Value = 7176.231
$ws.Range("L46").Value = 7176.231
$ws.Range("N46").Value = -7552.231
$ws.Range("H61").Value = 3064.7058
$ws.Range("I61").Value = 2034
$ws.Range("J61").Value = 5927.778
$ws.Range("K61").Value = 2034
$ws.Range("L61").Value = 5927.778
$ws.Range("M61").Value = -1832
$ws.Range("N61").Value = -6331.778
$ws.Range("H93").Value = 21962.375
$ws.Range("I93").Value = 15925
$ws.Range("J93").Value = 27999.75
$ws.Range("K93").Value = 15925
$ws.Range("L93").Value = 27999.75
$ws.Range("M93").Value = -14677
$ws.Range("N93").Value = -30495.75
$ws.Range("H100").Value = 7349.8335
$ws.Range("I100").Value = 3112.5
$ws.Range("J100").Value = 15824.5
$ws.Range("K100").Value = 3112.5
$ws.Range("L100").Value = 15824.5
$ws.Range("M100").Value = -2571.5
$ws.Range("N100").Value = -16906.5
$ws.Range("H109").Value = 17666.666
$ws.Range("J109").Value = 17666.666
$ws.Range("L109").Value = 17666.666
$ws.Range("N109").Value = -20440.666
$ws.Range("H113").Value = 3064.7058
$ws.Range("I113").Value = 2034
$ws.Range("J113").Value = 5927.778
$ws.Range("K113").Value = 2034
$ws.Range("L113").Value = 5927.778
$ws.Range("M113").Value = 136
$ws.Range("N113").Value = -10267.778

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2492.25
$ws.Range("I17").Value = 2495
$ws.Range("K17").Value = 2495
$ws.Range("M17").Value = -2323
$ws.Range("H62").Value = 3184.6155
$ws.Range("I62").Value = 2870
$ws.Range("K62").Value = 2870
$ws.Range("M62").Value = -2246
$ws.Range("H65").Value = 3184.6155
$ws.Range("I65").Value = 2870
$ws.Range("K65").Value = 14350
$ws.Range("M65").Value = -11230
$ws.Range("H81").Value = 699
$ws.Range("I81").Value = 699
$ws.Range("K81").Value = 1398
$ws.Range("M81").Value = -337
$ws.Range("H84").Value = 699
$ws.Range("I84").Value = 699
$ws.Range("K84").Value = 6990
$ws.Range("M84").Value = -1686
$ws.Range("H113").Value = 2352.6572
$ws.Range("I113").Value = 2594.36
$ws.Range("K113").Value = 7783.08
$ws.Range("M113").Value = -5613.08
$ws.Range("H122").Value = 6397.4346
$ws.Range("I122").Value = 2974.3845
$ws.Range("J122").Value = 10847.4
$ws.Range("K122").Value = 8923.1535
$ws.Range("L122").Value = 32542.2
$ws.Range("M122").Value = -6473.1535
$ws.Range("N122").Value = -37442.2
$ws.Range("H138").Value = 100998.8
$ws.Range("J138").Value = 110624.25
$ws.Range("L138").Value = 110624.25
$ws.Range("N138").Value = -120904.25
